# Generate Report for Handback
# Update "last generated" timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for b7d12afe... row
$wsOverview.Range("G4").Value = "2016-09-02 04:48:29"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for b7d12afe... row
$wsZhCn.Range("H4").Value = "2016-09-02 04:48:25"
$wsZhCn.Range("K4").Value = "2016-09-02 04:48:43"

# de-de sheet: Correspond Handback DateTime for b7d12afe... row
$wsDeDe.Range("K4").Value = "2016-09-02 04:48:50"
